$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "西部材料"
$ws.Range("B2").Value = "平潭发展"
$ws.Range("C2").Value = "平潭发展"
$ws.Range("A3").Value = "航天发展"
$ws.Range("B3").Value = "山子高科"
$ws.Range("C3").Value = "雪人集团"
$ws.Range("A4").Value = "浙江世宝"
$ws.Range("B4").Value = "西部材料"
$ws.Range("C4").Value = "东百集团"
$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "雪人集团"
$ws.Range("C5").Value = "航天发展"
$ws.Range("A6").Value = "雪人集团"
$ws.Range("B6").Value = "东百集团"
$ws.Range("C6").Value = "永辉超市"
$ws.Range("A7").Value = "东百集团"
$ws.Range("B7").Value = "浙江世宝"
$ws.Range("C7").Value = "西部材料"
$ws.Range("A8").Value = "神剑股份"
$ws.Range("B8").Value = "永辉超市"
$ws.Range("C8").Value = "浙江世宝"
$ws.Range("A9").Value = "山子高科"
$ws.Range("B9").Value = "航天发展"
$ws.Range("C9").Value = "山子高科"
$ws.Range("A10").Value = "永辉超市"
$ws.Range("B10").Value = "航天电子"
$ws.Range("C10").Value = "四川长虹"
$ws.Range("A11").Value = "航天电子"
$ws.Range("B11").Value = "通鼎互联"
$ws.Range("C11").Value = "航天电子"
$ws.Range("A12").Value = "通鼎互联"
$ws.Range("B12").Value = "王子新材"
$ws.Range("C12").Value = "再升科技"
$ws.Range("A13").Value = "再升科技"
$ws.Range("B13").Value = "百利电气"
$ws.Range("C13").Value = "神剑股份"
$ws.Range("A14").Value = "顺灏股份"
$ws.Range("B14").Value = "海南发展"
$ws.Range("C14").Value = "通鼎互联"
$ws.Range("A15").Value = "王子新材"
$ws.Range("B15").Value = "神剑股份"
$ws.Range("C15").Value = "金圆股份"
$ws.Range("A16").Value = "海南发展"
$ws.Range("B16").Value = "中国中免"
$ws.Range("C16").Value = "顺灏股份"
$ws.Range("A17").Value = "金圆股份"
$ws.Range("B17").Value = "再升科技"
$ws.Range("C17").Value = "王子新材"
$ws.Range("A18").Value = "伯特利"
$ws.Range("B18").Value = "中国卫星"
$ws.Range("C18").Value = "百大集团"
$ws.Range("A19").Value = "中国卫星"
$ws.Range("B19").Value = "海南瑞泽"
$ws.Range("C19").Value = "中国神华"
$ws.Range("A20").Value = "久之洋"
$ws.Range("B20").Value = "上海九百"
$ws.Range("C20").Value = "合富中国"
$ws.Range("A21").Value = "中国中免"
$ws.Range("B21").Value = "金圆股份"
$ws.Range("C21").Value = "九牧王"
